$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '95.102.85'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.37%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.571.16'
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '657.26'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.93%  '
$ws.Range("E7").Value = '  -0.60%  '
$ws.Range("E8").Value = '  -1.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.999'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.569.19'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.83%  '
$ws.Range("E12").Value = '  +1.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.37'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.79%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.47'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.237.34'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.76%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '94.980.85'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.39%  '
$ws.Range("E17").Value = '  -0.56%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.574.61'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.79%  '
$ws.Range("E19").Value = '  -7.51%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.63'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.79'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.77%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.45'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '508.26'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.478'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.80'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.98%  '
$ws.Range("E26").Value = '  -1.92%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '95.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.68'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.72%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.763.17'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.76%  '
$ws.Range("E30").Value = '  -1.33%  '
$ws.Range("E31").Value = '  -0.69%  '
$ws.Range("E32").Value = '  -0.50%  '
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.05%  '
$ws.Range("E35").Value = '  -2.84%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.82'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.62%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.68'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +12.83%  '
$ws.Range("E38").Value = '  -2.32%  '
$ws.Range("E39").Value = '  +7.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '582.74'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.40%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("E42").Value = '  -0.86%  '
$ws.Range("E43").Value = '  -2.23%  '
$ws.Range("E44").Value = '  +3.70%  '
$ws.Range("E45").Value = '  +1.33%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '34.44'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +31.79%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.27'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.51%  '
$ws.Range("E48").Value = '  -1.70%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0414'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.59'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.08%  '
$ws.Range("E51").Value = '  +0.45%  '
